$wb = $excel.ActiveWorkbook

# Add the new "EmpList" worksheet after the last sheet ("User")
$userSheet = $wb.Worksheets.Item("User")
$newSheet = $wb.Worksheets.Add($null, $userSheet)
$newSheet.Name = "EmpList"

# Populate header row, mirroring the Admin/admin123 login columns plus a new name
$newSheet.Range("A1").Value = "Admin"
$newSheet.Range("B1").Value = "admin123"
$newSheet.Range("C1").Value = "Lisa"

# Selection on the new sheet
$newSheet.Range("B6").Select()

# Make the new sheet the active/selected tab, and the Emp sheet no longer selected
$newSheet.Activate()

$wb.Save()
